# -----------------------------------------------------------------------
# "Dodanie podziału treningu na części"
# Adds a "Trening" (training-phase) column F, replaces the sample rows with
# a full GPS session split into "Duza Gra" / "Mala Gra" phases, and switches
# column A from text timestamps to real Excel datetimes.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the training-phase column - reuse the same header
# formatting (bold font + border, centered) as the other header cells.
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Trening"

# Row 2
$ws.Range("A2").Value = 45686
$ws.Range("B2").Value = "'"    # blank (no Seconds for the session-start marker)
$ws.Range("C2").Value = "'"    # blank (no Velocity for the session-start marker)
$ws.Range("D2").Value = "'"    # blank (no Acceleration_SMA for the session-start marker)
$ws.Range("E2").Value = "10-15"
$ws.Range("F2").Value = "Duża Gra"

# Row 3
$ws.Range("A3").Value = 45686.47614664352
$ws.Range("B3").Value = 1186.9
$ws.Range("C3").Value = 5.01
$ws.Range("D3").Value = 0.8203748975481303
$ws.Range("E3").Value = "5-10"
$ws.Range("F3").Value = "Duża Gra"

# Row 4
$ws.Range("A4").Value = 45686.47821030093
$ws.Range("B4").Value = 1365.2
$ws.Range("C4").Value = 7.57
$ws.Range("D4").Value = 0.4794053392750878
$ws.Range("E4").Value = "5-10"
$ws.Range("F4").Value = "Duża Gra"

# Row 5
$ws.Range("A5").Value = 45686.47891631944
$ws.Range("B5").Value = 1426.2
$ws.Range("C5").Value = 7.03
$ws.Range("D5").Value = 0.5369732209614344
$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Duża Gra"

# Row 6
$ws.Range("A6").Value = 45686.48759224537
$ws.Range("B6").Value = 2175.8
$ws.Range("C6").Value = 10.84
$ws.Range("D6").Value = 2.528612136840823
$ws.Range("E6").Value = "10-15"
$ws.Range("F6").Value = "Mała Gra"

# Row 7
$ws.Range("A7").Value = 45686.48794293981
$ws.Range("B7").Value = 2206.1
$ws.Range("C7").Value = 12.01
$ws.Range("D7").Value = 2.773899997983659
$ws.Range("E7").Value = "10-15"
$ws.Range("F7").Value = "Mała Gra"

# Row 8
$ws.Range("A8").Value = 45686.49236886574
$ws.Range("B8").Value = 2588.5
$ws.Range("C8").Value = 11.5
$ws.Range("D8").Value = 2.638877425874985
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

# Row 9
$ws.Range("A9").Value = 45686.487940625
$ws.Range("B9").Value = 2205.9
$ws.Range("C9").Value = 9.82
$ws.Range("D9").Value = 2.577154687472753
$ws.Range("E9").Value = "5-10"
$ws.Range("F9").Value = "Mała Gra"

# Row 10
$ws.Range("A10").Value = 45686.49347650463
$ws.Range("B10").Value = 2684.2
$ws.Range("C10").Value = 9.63
$ws.Range("D10").Value = 2.581018277576991
$ws.Range("E10").Value = "5-10"
$ws.Range("F10").Value = "Mała Gra"

# Row 11
$ws.Range("A11").Value = 45686.49419409722
$ws.Range("B11").Value = 2746.2
$ws.Range("C11").Value = 9.95
$ws.Range("D11").Value = 2.511234283447266
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"

# The three blank cells above were entered with a leading apostrophe so they
# land as empty *text*, not as cleared/blank cells; now drop that text-quote
# formatting again so the cells keep the workbook default style.
$ws.Range("B2:D2").Style = "Normal"

# Re-format column A as real date/time values (was plain text before).
# Apply the lower-case variant first, then switch to the upper-case variant
# actually used by the column, matching the source workbook's two number
# formats (164 defined-but-unused, 165 applied).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

